$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 8748
$ws.Cells.Item(18, 9).Value = 4996.6665
$ws.Cells.Item(18, 10).Value = 20002
$ws.Cells.Item(18, 11).Value = 4996.6665
$ws.Cells.Item(18, 12).Value = 20002
$ws.Cells.Item(18, 13).Value = -4712.6665
$ws.Cells.Item(18, 14).Value = -20570
$ws.Cells.Item(28, 8).Value = 765.3333
$ws.Cells.Item(28, 9).Value = 418.4
$ws.Cells.Item(28, 11).Value = 418.4
$ws.Cells.Item(28, 13).Value = 66.60000000000002
$ws.Cells.Item(55, 8).Value = 126.9
$ws.Cells.Item(55, 9).Value = 133.44444
$ws.Cells.Item(55, 10).Value = 68
$ws.Cells.Item(55, 11).Value = 133.44444
$ws.Cells.Item(55, 12).Value = 68
$ws.Cells.Item(55, 13).Value = 80.55556000000001
$ws.Cells.Item(55, 14).Value = -496
$ws.Cells.Item(70, 8).Value = 4140
$ws.Cells.Item(70, 9).Value = 2200
$ws.Cells.Item(70, 10).Value = 4625
$ws.Cells.Item(70, 11).Value = 6600
$ws.Cells.Item(70, 12).Value = 13875
$ws.Cells.Item(70, 13).Value = -6330
$ws.Cells.Item(70, 14).Value = -14415
$ws.Cells.Item(73, 8).Value = 4140
$ws.Cells.Item(73, 9).Value = 2200
$ws.Cells.Item(73, 10).Value = 4625
$ws.Cells.Item(73, 11).Value = 6600
$ws.Cells.Item(73, 12).Value = 13875
$ws.Cells.Item(73, 13).Value = -5664
$ws.Cells.Item(73, 14).Value = -15747
$ws.Cells.Item(74, 8).Value = 170000
$ws.Cells.Item(74, 9).Value = 5000
$ws.Cells.Item(74, 11).Value = 5000
$ws.Cells.Item(74, 13).Value = -4064
$ws.Cells.Item(76, 8).Value = 5119.8
$ws.Cells.Item(76, 9).Value = 5119.8
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 11).Value = 5119.8
$ws.Cells.Item(76, 12).Value = 0
$ws.Cells.Item(76, 13).Value = -4804.8
$ws.Cells.Item(76, 14).ClearContents()
$ws.Cells.Item(77, 8).Value = 170000
$ws.Cells.Item(77, 9).Value = 5000
$ws.Cells.Item(77, 11).Value = 25000
$ws.Cells.Item(77, 13).Value = -20320
$ws.Cells.Item(79, 8).Value = 5119.8
$ws.Cells.Item(79, 9).Value = 5119.8
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 11).Value = 5119.8
$ws.Cells.Item(79, 12).Value = 0
$ws.Cells.Item(79, 13).Value = -4027.8
$ws.Cells.Item(79, 14).ClearContents()
$ws.Cells.Item(86, 8).Value = 11042.444
$ws.Cells.Item(86, 9).Value = 4197.4287
$ws.Cells.Item(86, 11).Value = 4197.4287
$ws.Cells.Item(86, 13).Value = -3074.4287
$ws.Cells.Item(89, 8).Value = 11042.444
$ws.Cells.Item(89, 9).Value = 4197.4287
$ws.Cells.Item(89, 11).Value = 20987.1435
$ws.Cells.Item(89, 13).Value = -15371.1435
$ws.Cells.Item(92, 8).Value = 587.4
$ws.Cells.Item(92, 9).Value = 641.6667
$ws.Cells.Item(92, 10).Value = 99
$ws.Cells.Item(92, 11).Value = 641.6667
$ws.Cells.Item(92, 12).Value = 99
$ws.Cells.Item(92, 13).Value = 606.3333
$ws.Cells.Item(92, 14).Value = -2595
$ws.Cells.Item(98, 8).Value = 1144.1
$ws.Cells.Item(98, 9).Value = 1144.1
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 1144.1
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 13).Value = 353.9000000000001
$ws.Cells.Item(98, 14).ClearContents()
$ws.Cells.Item(106, 8).Value = 2799.6667
$ws.Cells.Item(106, 9).Value = 1000
$ws.Cells.Item(106, 10).Value = 3699.5
$ws.Cells.Item(106, 11).Value = 1000
$ws.Cells.Item(106, 12).Value = 3699.5
$ws.Cells.Item(106, 13).Value = -369
$ws.Cells.Item(106, 14).Value = -4961.5
$ws.Cells.Item(113, 8).Value = 4716.75
$ws.Cells.Item(113, 9).Value = 3900
$ws.Cells.Item(113, 10).Value = 4989
$ws.Cells.Item(113, 11).Value = 3900
$ws.Cells.Item(113, 12).Value = 4989
$ws.Cells.Item(113, 13).Value = -646
$ws.Cells.Item(113, 14).Value = -11497
$ws.Cells.Item(116, 8).Value = 5833.1665
$ws.Cells.Item(116, 9).Value = 4749.75
$ws.Cells.Item(116, 11).Value = 4749.75
$ws.Cells.Item(116, 13).Value = -1307.75
$ws.Cells.Item(122, 8).Value = 1144.1
$ws.Cells.Item(122, 9).Value = 1144.1
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 3432.3
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -982.2999999999997
$ws.Cells.Item(122, 14).ClearContents()

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 13).ClearContents()
$ws.Cells.Item(32, 8).Value = 8866.825999999999
$ws.Cells.Item(32, 9).Value = 7613.9414
$ws.Cells.Item(32, 11).Value = 7613.9414
$ws.Cells.Item(32, 13).Value = -7326.9414
$ws.Cells.Item(45, 8).Value = 1959.8
$ws.Cells.Item(45, 9).Value = 1959.8
$ws.Cells.Item(45, 11).Value = 1959.8
$ws.Cells.Item(45, 13).Value = -1582.8
$ws.Cells.Item(61, 8).Value = 8001.25
$ws.Cells.Item(61, 9).Value = 7715.7144
$ws.Cells.Item(61, 10).Value = 10000
$ws.Cells.Item(61, 11).Value = 7715.7144
$ws.Cells.Item(61, 12).Value = 10000
$ws.Cells.Item(61, 13).Value = -7503.7144
$ws.Cells.Item(61, 14).Value = -10424
$ws.Cells.Item(110, 8).Value = 3999.5
$ws.Cells.Item(110, 9).Value = 3999.5
$ws.Cells.Item(110, 11).Value = 3999.5
$ws.Cells.Item(110, 13).Value = -1954.5
$ws.Cells.Item(116, 8).Value = 0
$ws.Cells.Item(116, 9).Value = 0
$ws.Cells.Item(116, 11).Value = 0
$ws.Cells.Item(116, 13).ClearContents()
$ws.Cells.Item(136, 8).Value = 8001.25
$ws.Cells.Item(136, 9).Value = 7715.7144
$ws.Cells.Item(136, 10).Value = 10000
$ws.Cells.Item(136, 11).Value = 23147.1432
$ws.Cells.Item(136, 12).Value = 30000
$ws.Cells.Item(136, 13).Value = -20597.1432
$ws.Cells.Item(136, 14).Value = -35100

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 13).ClearContents()
$ws.Cells.Item(94, 8).Value = 2384.7778
$ws.Cells.Item(94, 9).Value = 2384.7778
$ws.Cells.Item(94, 11).Value = 2384.7778
$ws.Cells.Item(94, 13).Value = -1933.7778

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 131.4375
$ws.Cells.Item(7, 9).Value = 77.23077000000001
$ws.Cells.Item(7, 10).Value = 366.33334
$ws.Cells.Item(7, 11).Value = 77.23077000000001
$ws.Cells.Item(7, 12).Value = 366.33334
$ws.Cells.Item(7, 13).Value = 35.76922999999999
$ws.Cells.Item(7, 14).Value = -592.33334
$ws.Cells.Item(31, 8).Value = 5361.3076
$ws.Cells.Item(31, 9).Value = 4349.625
$ws.Cells.Item(31, 11).Value = 4349.625
$ws.Cells.Item(31, 13).Value = -4054.625
$ws.Cells.Item(34, 8).Value = 5361.3076
$ws.Cells.Item(34, 9).Value = 4349.625
$ws.Cells.Item(34, 11).Value = 4349.625
$ws.Cells.Item(34, 13).Value = -4147.625
$ws.Cells.Item(58, 8).Value = 5718.5
$ws.Cells.Item(58, 9).Value = 1437.3334
$ws.Cells.Item(58, 11).Value = 1437.3334
$ws.Cells.Item(58, 13).Value = -1234.3334
$ws.Cells.Item(68, 8).Value = 30439.533
$ws.Cells.Item(68, 9).Value = 19513.285
$ws.Cells.Item(68, 11).Value = 19513.285
$ws.Cells.Item(68, 13).Value = -18764.285
$ws.Cells.Item(71, 8).Value = 30439.533
$ws.Cells.Item(71, 9).Value = 19513.285
$ws.Cells.Item(71, 11).Value = 58539.855
$ws.Cells.Item(71, 13).Value = -54795.855
$ws.Cells.Item(81, 8).Value = 67888.5
$ws.Cells.Item(81, 9).Value = 0
$ws.Cells.Item(81, 11).Value = 0
$ws.Cells.Item(81, 13).ClearContents()
$ws.Cells.Item(82, 8).Value = 0
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 14).ClearContents()
$ws.Cells.Item(84, 8).Value = 67888.5
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 13).ClearContents()
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 14).ClearContents()
$ws.Cells.Item(107, 8).Value = 1798
$ws.Cells.Item(107, 9).Value = 2497.5
$ws.Cells.Item(107, 11).Value = 2497.5
$ws.Cells.Item(107, 13).Value = -577.5
$ws.Cells.Item(136, 8).Value = 5718.5
$ws.Cells.Item(136, 9).Value = 1437.3334
$ws.Cells.Item(136, 11).Value = 4312.0002
$ws.Cells.Item(136, 13).Value = -1762.0002

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(51, 8).Value = 1003.4286
$ws.Cells.Item(51, 10).Value = 1005
$ws.Cells.Item(51, 12).Value = 3015
$ws.Cells.Item(51, 14).Value = -3935

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(136, 8).Value = 0
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 14).ClearContents()

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(75, 8).Value = 60000
$ws.Cells.Item(75, 10).Value = 60000
$ws.Cells.Item(75, 12).Value = 60000
$ws.Cells.Item(75, 14).Value = -61872
$ws.Cells.Item(78, 8).Value = 60000
$ws.Cells.Item(78, 10).Value = 60000
$ws.Cells.Item(78, 12).Value = 180000
$ws.Cells.Item(78, 14).Value = -189360
$ws.Cells.Item(123, 8).Value = 79996
$ws.Cells.Item(123, 10).Value = 79996
$ws.Cells.Item(123, 12).Value = 79996
$ws.Cells.Item(123, 14).Value = -89796
$ws.Cells.Item(132, 8).Value = 13142.714
$ws.Cells.Item(132, 9).Value = 5999.75
$ws.Cells.Item(132, 11).Value = 17999.25
$ws.Cells.Item(132, 13).Value = -15469.25

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(75, 8).Value = 24000
$ws.Cells.Item(75, 10).Value = 24000
$ws.Cells.Item(75, 12).Value = 24000
$ws.Cells.Item(75, 14).Value = -25872
$ws.Cells.Item(78, 8).Value = 24000
$ws.Cells.Item(78, 10).Value = 24000
$ws.Cells.Item(78, 12).Value = 72000
$ws.Cells.Item(78, 14).Value = -81360
$ws.Cells.Item(132, 8).Value = 2232.7778
$ws.Cells.Item(132, 9).Value = 2199.375
$ws.Cells.Item(132, 10).Value = 2500
$ws.Cells.Item(132, 11).Value = 6598.125
$ws.Cells.Item(132, 12).Value = 7500
$ws.Cells.Item(132, 13).Value = -4068.125
$ws.Cells.Item(132, 14).Value = -12560
